$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for D-column price cells whose values look numeric,
# so Excel stores them as text (matching the source inlineStr type) instead
# of silently converting to a Number and dropping trailing zeros.
$ws.Range("D2").Value = '67.229.27'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '2.551.07'
$ws.Range("E3").Value = '  -2.27%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.81'
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.67'
$ws.Range("E6").Value = '  +5.11%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.530'
$ws.Range("E8").Value = '  +0.58%  '
$ws.Range("D9").Value = '2.550.08'
$ws.Range("E9").Value = '  -2.37%  '
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("E11").Value = '  +1.83%  '
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("E13").Value = '  -4.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.96'
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").Value = '3.011.66'
$ws.Range("E15").Value = '  -2.48%  '
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").Value = '67.104.36'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").Value = '2.557.51'
$ws.Range("E18").Value = '  -2.67%  '
$ws.Range("E19").Value = '  +3.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.33'
$ws.Range("E20").Value = '  -2.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '355.78'
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.20'
$ws.Range("E22").Value = '  -1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.65'
$ws.Range("E23").Value = '  +0.96%  '
$ws.Range("E24").Value = '  +6.41%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.94'
$ws.Range("E26").Value = '  +1.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.07'
$ws.Range("E27").Value = '  -3.82%  '
$ws.Range("D28").Value = '2.673.06'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '533.40'
$ws.Range("E31").Value = '  -0.95%  '
$ws.Range("E32").Value = '  +0.92%  '
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("E35").Value = '  -0.41%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  +0.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.64'
$ws.Range("E38").Value = '  -1.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.70'
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("E40").Value = '  +1.20%  '
$ws.Range("E41").Value = '  -1.52%  '
$ws.Range("E42").Value = '  +0.36%  '
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("E44").Value = '  +5.10%  '
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '150.07'
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("E48").Value = '  -1.90%  '
$ws.Range("E49").Value = '  -4.95%  '
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("E51").Value = '  +0.57%  '
